# chore: adapt column header formatting to respective input file names
#
# Renames the "_old" / "_new" header-name suffixes used throughout the
# AHB-Diff header row to the concrete format-version identifiers
# "_FV2410" (old/previous format version) and "_FV2504" (new/current
# format version), then turns the sheet's data range into a proper
# Excel Table ("Table1") and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the column headers in row 1 --------------------------------
# Columns A:J carried the "_old" suffix, columns L:U carried the "_new"
# suffix. Column K ("diff") is untouched.
$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

$fv2504Headers = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

# Columns A-J (1-10) -> "_old" becomes "_FV2410"
for ($i = 0; $i -lt $fv2410Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2410Headers[$i]
}

# Column K (11) stays "diff" - nothing to do.

# Columns L-U (12-21) -> "_new" becomes "_FV2504"
for ($i = 0; $i -lt $fv2504Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2504Headers[$i]
}

# --- 2) Convert the used range into an Excel Table ("Table1") -------------
$dataRange = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# --- 3) Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
